$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 65: append two more lines to the existing comment, and bump the hours ---
$existingG65 = $ws.Cells.Item(65, 7).Value2
$newG65 = $existingG65 + "`nBug gefixed, bei dem Auswählen des default download directories nicht mehr möglich war (Fehler im Umgang mit preferences)`nreadme.md aktualisiert"
$ws.Cells.Item(65, 7).Value2 = $newG65
$ws.Cells.Item(65, 6).Value2 = 5

# --- New row 66: copy formatting from row 65 (date / hours / comment columns) ---
$ws.Range("E65:G65").Copy() | Out-Null
$ws.Range("E66:G66").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Cells.Item(66, 5).Value2 = 43840
$ws.Cells.Item(66, 6).Value2 = 4
$ws.Cells.Item(66, 7).Value2 = "Codedokumentierung`nAnpassungen für TdOT (Einfügen neuer Seite, Counter bis Release)`nDiverse Vorbereitungen für TdOT getroffen"

# Match the row heights recorded in the saved workbook (Excel auto-sized these
# to fit the wrapped, multi-line text).
$ws.Rows.Item(65).RowHeight = 72
$ws.Rows.Item(66).RowHeight = 43.2

# Recalculate so C5 (=SUM(F:F)) reflects the updated/added hours.
$excel.CalculateFull() | Out-Null

# --- Update the view: scroll down a bit and select G66:G68 ---
$ws.Range("G66:G68").Select() | Out-Null
